$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "49,999"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "467,701"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "1,303,395"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "784,464"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.16"

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "1,978,838"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "1,366,330"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "32.14"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "50,000"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "600,000"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "11"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1,250,024,999"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "1,250,024,999"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "196.94"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "49,999"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "625,025,000"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1198.34"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "965,541"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "28,023"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "1,355,853"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "784,464"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "1,366,047"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "698,893"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "49.00"

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "324,631,638"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "324,597,122"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "444.99"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "49,999"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "400,055"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "290.92"

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "1,115,742"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "89,809"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.45"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "1,645,530"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "784,464"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.48"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "1,813,823"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "1,179,702"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "46.59"

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "621,490,726"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "621,490,273"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "849.26"

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "49,999"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "258,931"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "291.75"

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "1,120,072"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "93,497"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.80"

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "2,300,148"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "784,464"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.53"

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "1,419,654"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "756,351"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.71"
